# ------------------------------------------------------------------
# fil_rouge.docx edit script
#
# 1. Append a new sentence (as its own run) after "... envoyés et reçus"
# 2. Re-flow "grâce à notre système ... accéder/accédé" text: fix the
#    "accéder" -> "accédé" typo, drop the lastRenderedPageBreak that
#    used to sit on the "système..." run, and split the text into
#    three runs with the "_GoBack" bookmark moved to the new seam.
# 3. Remove the (now relocated) "_GoBack" bookmark that used to sit
#    after "Cahier des Charges".
# 4. Remove the lastRenderedPageBreak that used to sit before "3 - ".
#
# NOTE: inserting/retyping text right next to existing runs normally
# gets silently merged back into neighbouring runs whenever they end
# up with identical formatting (this engine re-coalesces adjacent
# same-format runs within a paragraph after any edit). To keep seams
# as distinct <w:r> elements (as the target XML wants) we drop a
# throw-away bookmark exactly on the seam and delete it again
# immediately - its mere presence during save is enough to keep the
# runs from being re-joined, and removing it leaves no trace behind.
# All run-splitting bookmarks are therefore added only at the very
# end, after every text edit in a paragraph is already in place.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# Apostrophe (U+2019) used throughout the source document.
$apos = [char]0x2019

# ------------------------------------------------------------------
# 1) "Ils ne peuvent proposer ... et reçus" -> add trailing sentence
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Ils ne peuvent proposer un service Web permettant la lecture des messages envoyés et reçus", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamPt = $rng.End
$rng.InsertAfter(", ainsi que l" + $apos + "envoie de nouveaux messages.")
$seam = $d.Range($seamPt, $seamPt)
$d.Bookmarks.Add("zzSeam1", $seam) | Out-Null
$d.Bookmarks("zzSeam1").Delete()

# ------------------------------------------------------------------
# 2) "grâce à notre ... accéder par plusieurs plateformes" rework
# ------------------------------------------------------------------
# 2a. Fix the typo accéder -> accédé
$rng = $d.Content
$rng.Find.Execute("peut-être accéder par", $true, $false, $false, $false, $false, $true, 1, $false, `
    "peut-être accédé par", 2) | Out-Null

# 2b. Remove the lastRenderedPageBreak that sat at the start of the
#     run carrying "système, un compte ...". Cut out just its first
#     character and retype it - this drops the stale layout hint
#     without otherwise touching the paragraph text.
$rng = $d.Content
$rng.Find.Execute("système, un compte utilisateur peut-être accédé par", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart = $rng.Start
$firstChar = $d.Range($segStart, $segStart + 1)
$firstChar.Delete()
$reins = $d.Range($segStart, $segStart)
$reins.InsertAfter("s")

# 2c. Re-create the run seams (rightmost first so earlier Find calls
#     in this same pass are not disturbed by the zero-width bookmark
#     inserts):
#       ... abord.  |  En effet, grâce à notre système, un comp | te utilisateur peut-être accédé | <bookmark _GoBack> |  par plusieurs plateformes. ...
$rng = $d.Content
$rng.Find.Execute("système, un compte utilisateur peut-être accédé", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamGoBack = $rng.End
$seam = $d.Range($seamGoBack, $seamGoBack)
$d.Bookmarks.Add("_GoBack", $seam) | Out-Null

$rng = $d.Content
$rng.Find.Execute("système, un comp", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamComp = $rng.End
$seam = $d.Range($seamComp, $seamComp)
$d.Bookmarks.Add("zzSeam2", $seam) | Out-Null
$d.Bookmarks("zzSeam2").Delete()

$rng = $d.Content
$rng.Find.Execute("au premier abord.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamAbord = $rng.End
$seam = $d.Range($seamAbord, $seamAbord)
$d.Bookmarks.Add("zzSeam3", $seam) | Out-Null
$d.Bookmarks("zzSeam3").Delete()

$rng = $d.Content
$rng.Find.Execute("depuis telle ou telle plateforme.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamPlateforme = $rng.End
$seam = $d.Range($seamPlateforme, $seamPlateforme)
$d.Bookmarks.Add("zzSeam4", $seam) | Out-Null
$d.Bookmarks("zzSeam4").Delete()

# ------------------------------------------------------------------
# 4) Drop the lastRenderedPageBreak that sat before "3 - " (Titre1)
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("3 - La gestion du temps et du projet", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$segStart = $rng.Start
$firstChar = $d.Range($segStart, $segStart + 1)
$firstChar.Delete()
$reins = $d.Range($segStart, $segStart)
$reins.InsertAfter("3")

# Re-create the seam between "3 - " and "La gestion du temps et du projet"
$rng = $d.Content
$rng.Find.Execute("3 - ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$seamTitre = $rng.End
$seam = $d.Range($seamTitre, $seamTitre)
$d.Bookmarks.Add("zzSeam5", $seam) | Out-Null
$d.Bookmarks("zzSeam5").Delete()

Write-Output "done"
